{"js": "// Update the benchmark stats table in place.\n// The document contains a single, single-column table whose rows hold\n// one metric value each (plus three \"legacy\" rows near the bottom that\n// still carry the old multi-column/tab-separated layout from an earlier\n// version of the report). This edit:\n//   1. Refreshes a handful of metric values in the first block of rows.\n//   2. Collapses the three leftover multi-run/tab rows down to the single\n//      summary value they should have held all along.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> new cell text\nconst cellUpdates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"699\",\n  4: \"0.00002\",\n  6: \"0.00019\",\n  7: \"0.00009\",\n  8: \"0.00036\",\n  9: \"0.00042\",\n  10: \"0.00048\",\n  11: \"0.16135\",\n  // Legacy rows: collapse the tab-separated run list down to one value.\n  43: \"99.85\",\n  44: \"0.16\",\n  45: \"108\",\n};\n\nfor (const rowIndexStr of Object.keys(cellUpdates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newText = cellUpdates[rowIndexStr];\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.paragraphs.getFirst().getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the benchmark stats table in place.\n# The document contains a single, single-column table whose rows hold\n# one metric value each (plus three \"legacy\" rows near the bottom that\n# still carry the old multi-column/tab-separated layout from an earlier\n# version of the report). This edit:\n#   1. Refreshes a handful of metric values in the first block of rows.\n#   2. Collapses the three leftover multi-run/tab rows down to the single\n#      summary value they should have held all along.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Word COM table rows/columns are 1-based.\n$cellUpdates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"699\"\n    5  = \"0.00002\"\n    7  = \"0.00019\"\n    8  = \"0.00009\"\n    9  = \"0.00036\"\n    10 = \"0.00042\"\n    11 = \"0.00048\"\n    12 = \"0.16135\"\n    # Legacy rows: collapse the tab-separated run list down to one value.\n    44 = \"99.85\"\n    45 = \"0.16\"\n    46 = \"108\"\n}\n\nforeach ($rowIndex in $cellUpdates.Keys) {\n    $cell = $tbl.Cell($rowIndex, 1)\n    $cell.Range.Text = $cellUpdates[$rowIndex]\n}\n"}
